$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.735172271728516
$ws.Range("B1").Value = 2.678700923919678
$ws.Range("C1").Value = 2.039378643035889
$ws.Range("D1").Value = 1.952785611152649
$ws.Range("E1").Value = 2.086557626724243
